$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: "Lowest monthly average" label + MIN formula ---
$ws.Range("A14").Value = "Lowest monthly average"
$ws.Range("B14").Formula = "=MIN(B12:M12)"

# --- Row 15: MAX formula (no label in this row, matches the source sheet) ---
$ws.Range("B15").Formula = "=MAX(B12:M12)"

# --- Row 16: "Highest monthly average" label ---
$ws.Range("A16").Value = "Highest monthly average"
$ws.Rows.Item(16).RowHeight = 15.75

# --- Widen column A so the new, longer labels fit ---
$ws.Columns.Item(1).ColumnWidth = 20.8

# --- Update the window scroll position / selection to match the saved view ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
[void]$ws.Range("R14").Select()
